$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 137 (shifts existing rows 137:155 down to 138:156)
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row 137 with the new data record
$ws.Range("A137").Value = 7
$ws.Range("B137").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C137").Value = "Ñuble"
$ws.Range("D137").Value = 44449
$ws.Range("E137").Value = 16
$ws.Range("F137").Value = 100112008
$ws.Range("G137").Value = "Coliflor"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 300
$ws.Range("K137").Value = 700
$ws.Range("L137").Value = 750
$ws.Range("M137").Value = 725
$ws.Range("N137").Value = "$/unidad"
$ws.Range("O137").Value = "Región del Maule"
$ws.Range("P137").Value = 725
$ws.Range("Q137").Value = 1
$ws.Range("R137").Value = "Hortaliza"
